$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at the top; existing rows 1-12 shift down to 2-13
$ws.Rows.Item(1).EntireRow.Insert()

# Populate the new header row with the column labels
$ws.Range("A1").Value = "V"
$ws.Range("B1").Value = "F"
$ws.Range("C1").Value = "A"
$ws.Range("D1").Value = "Life"

# Match the saved selection state (active cell F2)
$ws.Range("F2").Select() | Out-Null
